$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-Str($s) {
    $s = $s.Replace("D80", "D86")
    $s = $s.Replace("D51", "D55")
    $s = $s.Replace("D64", "D69")
    $s = $s.Replace("S30", "S31")
    return $s
}

$ur = $ws.UsedRange
$rows = $ur.Rows.Count
$cols = $ur.Columns.Count

$changed = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newval = Transform-Str $val
            if ($newval -ne $val) {
                $cell.Value2 = $newval
                $changed = $changed + 1
            }
        }
    }
}

Write-Output "Cells changed: $changed"
